$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.356.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "'1.857.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.73%  "
$ws.Range("D5").Value = "'314.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").Value = "'0.4615"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("D8").Value = "'0.3714"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.07322"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").Value = "'0.8817"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").Value = "'19.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "'0.07807"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "'1.865.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").Value = "'5.383"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'6.548"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "'91.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "'0.000009071"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.24%  "
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "'14.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").Value = "'27.371.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "'2.101.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.02%  "
$ws.Range("D25").Value = "'1.933"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.54%  "
$ws.Range("D26").Value = "'152.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "'18.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("D28").Value = "'2.074"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").Value = "'5.111"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "'115.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "'0.08860"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "'0.7725"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.06%  "
$ws.Range("D33").Value = "'3.042"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.12%  "
$ws.Range("D34").Value = "'1.176"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.52%  "
$ws.Range("D35").Value = "'4.497"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").Value = "'2.651"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.91%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "'0.01958"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").Value = "'0.05230"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'2.952"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").Value = "'7.038"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("D42").Value = "'0.5148"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").Value = "'0.1636"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").Value = "'8.426"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.48%  "
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").Value = "'10.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").Value = "'103.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").Value = "'1.651"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("D50").Value = "'0.06225"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "'65.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.26%  "
